# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "33.599.48"
$ws.Range("E2").Value = "  -1.08%  "
# Row 3
$ws.Range("D3").Value = "1.763.59"
$ws.Range("E3").Value = "  -1.11%  "
# Row 4
$ws.Range("E4").Value = "  +0.18%  "
# Row 5
$ws.Range("D5").Value = "'222.97"
# Row 6
$ws.Range("D6").Value = "'0.542"
$ws.Range("E6").Value = "  -1.82%  "
# Row 7
$ws.Range("E7").Value = "  +0.22%  "
# Row 8
$ws.Range("D8").Value = "'31.85"
$ws.Range("E8").Value = "  +1.22%  "
# Row 9
$ws.Range("D9").Value = "'0.286"
$ws.Range("E9").Value = "  -0.40%  "
# Row 10
$ws.Range("D10").Value = "'0.0684"
$ws.Range("E10").Value = "  -3.64%  "
# Row 11
$ws.Range("D11").Value = "'0.0936"
$ws.Range("E11").Value = "  +1.67%  "
# Row 12
$ws.Range("D12").Value = "2.011.95"
$ws.Range("E12").Value = "  -1.39%  "
# Row 13
$ws.Range("D13").Value = "'11.06"
$ws.Range("E13").Value = "  +5.16%  "
# Row 14
$ws.Range("D14").Value = "1.774.28"
$ws.Range("E14").Value = "  -0.53%  "
# Row 15
$ws.Range("D15").Value = "33.597.68"
$ws.Range("E15").Value = "  -1.15%  "
# Row 16
$ws.Range("D16").Value = "'0.607"
$ws.Range("E16").Value = "  -3.07%  "
# Row 17
$ws.Range("D17").Value = "'4.10"
$ws.Range("E17").Value = "  -2.55%  "
# Row 18
$ws.Range("D18").Value = "'66.36"
$ws.Range("E18").Value = "  -2.38%  "
# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0769"
$ws.Range("E19").Value = "  -1.34%  "
# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'236.63"
$ws.Range("E20").Value = "  -3.31%  "
# Row 21
$ws.Range("E21").Value = "  +0.39%  "
# Row 22
$ws.Range("D22").Value = "'10.52"
$ws.Range("E22").Value = "  -1.62%  "
# Row 23
$ws.Range("D23").Value = "'4.00"
$ws.Range("E23").Value = "  -1.73%  "
# Row 24
$ws.Range("E24").Value = "  -2.90%  "
# Row 25
$ws.Range("D25").Value = "'158.92"
$ws.Range("E25").Value = "  +0.93%  "
# Row 26
$ws.Range("D26").Value = "'16.04"
$ws.Range("E26").Value = "  -2.08%  "
# Row 27
$ws.Range("D27").Value = "'6.98"
$ws.Range("E27").Value = "  -0.14%  "
# Row 28
$ws.Range("E28").Value = "  -0.78%  "
# Row 29
$ws.Range("E29").Value = "  +0.28%  "
# Row 30
$ws.Range("E30").Value = "  +0.91%  "
# Row 31
$ws.Range("D31").Value = "'0.0509"
$ws.Range("E31").Value = "  -2.41%  "
# Row 32
$ws.Range("D32").Value = "'3.58"
$ws.Range("E32").Value = "  -2.97%  "
# Row 33
$ws.Range("D33").Value = "'3.47"
$ws.Range("E33").Value = "  -0.60%  "
# Row 34
$ws.Range("D34").Value = "'1.77"
$ws.Range("E34").Value = "  -2.31%  "
# Row 35
$ws.Range("D35").Value = "1.374.34"
$ws.Range("E35").Value = "  -1.93%  "
# Row 36
$ws.Range("D36").Value = "'0.643"
$ws.Range("E36").Value = "  +0.51%  "
# Row 37
$ws.Range("E37").Value = "  -2.47%  "
# Row 38
$ws.Range("E38").Value = "  -1.59%  "
# Row 39
$ws.Range("E39").Value = "  +1.14%  "
# Row 40
$ws.Range("D40").Value = "'2.20"
$ws.Range("E40").Value = "  +4.73%  "
# Row 41
$ws.Range("D41").Value = "'77.31"
$ws.Range("E41").Value = "  -2.71%  "
# Row 42
$ws.Range("D42").Value = "'2.65"
$ws.Range("E42").Value = "  -2.43%  "
# Row 43
$ws.Range("D43").Value = "'0.897"
$ws.Range("E43").Value = "  -4.02%  "
# Row 44
$ws.Range("D44").Value = "'13.33"
$ws.Range("E44").Value = "  +12.98%  "
# Row 45
$ws.Range("E45").Value = "  +4.44%  "
# Row 46
$ws.Range("D46").Value = "0.0₆0136"
$ws.Range("E46").Value = "  +14.42%  "
# Row 47
$ws.Range("D47").Value = "'0.0497"
$ws.Range("E47").Value = "  +1.23%  "
# Row 48
$ws.Range("D48").Value = "'106.44"
$ws.Range("E48").Value = "  +0.86%  "
# Row 49
$ws.Range("D49").Value = "'5.78"
$ws.Range("E49").Value = "  -2.65%  "
# Row 50
$ws.Range("D50").Value = "1.914.39"
$ws.Range("E50").Value = "  -0.79%  "
# Row 51
$ws.Range("E51").Value = "  +0.57%  "
